$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting existing rows (2..16) down to (3..17)
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the "Not applicable" entry
$ws.Range("A2").Value = -1
$ws.Range("B2").Value = "Not applicable"

# Update the named range to cover the extra row
$wb.Names.Item("dbo_biogroup").RefersTo = "=dbo_biogroup!`$A`$1:`$B`$14"

# Move the active selection to D3 (matches the post-edit selection state)
$ws.Range("D3").Select()
